# Add a new "Modelo" column (F) describing the model used, matching the
# values already present for each disease row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n" +
             "                                            estimator=Pipeline(steps=[('model',`n" +
             "                                                                       LinearRegression())]),`n" +
             "                                            param_grid={'model__fit_intercept': [True,`n" +
             "                                                                                 False]},`n" +
             "                                            scoring='neg_mean_squared_error'))"

# Header cell F1, formatted like the other header cells (A1:E1).
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells F2:F10 with the model description.
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 6).Value = $modelText
}

# The multi-line text above can make Excel auto-expand the row heights;
# restore them to the sheet's default (no explicit/custom row height).
$ws.Range("A2:A10").EntireRow.AutoFit()
